$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item("magapoke_2025-10-29")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "magapoke_2025-11-05"

$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$srcSheet.Range("A1:B1").Copy($ws.Range("A1:B1"))

$titles = @(
  '黒月のイェルクナハト',
  'ドリーム☆ジャンボ☆ガール',
  'K-9~警視庁公安部公安第9課異能対策係~',
  'アイドラトリィ',
  '黄昏町プリズナーズ',
  'せいぶつ部の田辺くん',
  '春くらり',
  '篝家の８兄弟',
  'ハードワーカー中田',
  '生きたがりの人狼',
  'ナキナギ',
  'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
  'スルガメテオ',
  '夜鐘のキト',
  '屋根の下のアルテミス',
  'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
  '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
  '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
  'MYS',
  'ハナバス　苔石花江のバスケ論',
  '卒業アルバムの彼女たち',
  'それがメイドのカンナです',
  '英雄と魔女の転生ラブコメ',
  'ともだちづくり',
  '明智ナンバーワン',
  'ナマイキ旭ちゃんをわからせたい',
  '追放されなかった男　～二度目の人生は土下座から始まりました～',
  '永久のユウグレ',
  'じゅーくぼっくす',
  '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
  'JK Biker',
  '平成転生',
  '鳴るさんだぁ',
  'ハプスブルク家の華麗なる受難',
  '人生逆転ダンジョン',
  '眠れる森のレガ',
  '花子狩り',
  '〈小市民〉 春期限定いちごタルト事件',
  '鉱石令嬢〜没落した悪役令嬢が炭鉱で一山当てるまでのお話〜',
  '東京デスレース',
  '白銀のキュイジーヌ～明治外交官の料理人～',
  'イエティ、とある日々'
)

for ($i = 0; $i -lt $titles.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $i + 1
  $ws.Cells.Item($row, 2).Value = $titles[$i]
}
